$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 110, shifting existing rows 110:208 down to 111:209
$ws.Rows.Item(110).Insert()

# Populate the newly inserted row with its data
$ws.Cells.Item(110, 1).Value = 8
$ws.Cells.Item(110, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(110, 3).Value = "Coquimbo"
$ws.Cells.Item(110, 4).Value = 44512
$ws.Cells.Item(110, 5).Value = 4
$ws.Cells.Item(110, 6).Value = 100112032
$ws.Cells.Item(110, 7).Value = "Zapallo italiano"
$ws.Cells.Item(110, 8).Value = "Sin especificar"
$ws.Cells.Item(110, 9).Value = "Primera"
$ws.Cells.Item(110, 10).Value = 520
$ws.Cells.Item(110, 11).Value = 10000
$ws.Cells.Item(110, 12).Value = 11000
$ws.Cells.Item(110, 13).Value = 10500
$ws.Cells.Item(110, 14).Value = "`$/caja 70 unidades"
$ws.Cells.Item(110, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(110, 16).Value = 150
$ws.Cells.Item(110, 17).Value = 70
$ws.Cells.Item(110, 18).Value = "Hortaliza"
